# Update kim+sam example for sharing
$wb = $excel.ActiveWorkbook

# Update the "Fixed Assets" sheet (the active sheet) with a new row describing
# the house / residence fixed asset.
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = $true
$ws.Range("B2").Value = "house"
$ws.Range("C2").Value = "residence"
$ws.Range("D2").Value = 2020
$ws.Range("E2").Value = 400000
$ws.Range("F2").Value = 600000
$ws.Range("G2").Value = 3.6
$ws.Range("H2").Value = 2090
$ws.Range("I2").Value = 5

# Move / update the selected cell to A3, matching the saved workbook state.
$ws.Range("A3").Select() | Out-Null
